$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.959.82'
$ws.Range('E2').Value = '  +2.84%  '
$ws.Range('D3').Value = '1.600.33'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '212.44'
$ws.Range('E5').Value = '  +2.83%  '
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('E7').Value = '  +1.19%  '
$ws.Range('D8').Value = '0.247'
$ws.Range('E8').Value = '  +1.98%  '
$ws.Range('E9').Value = '  +0.76%  '
$ws.Range('D10').Value = '18.08'
$ws.Range('E10').Value = '  +2.04%  '
$ws.Range('E11').Value = '  +4.11%  '
$ws.Range('D12').Value = '1.823.11'
$ws.Range('E12').Value = '  +2.97%  '
$ws.Range('D13').Value = '1.593.21'
$ws.Range('E13').Value = '  +2.53%  '
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').Value = '0.509'
$ws.Range('E15').Value = '  +0.95%  '
$ws.Range('D16').Value = '25.974.27'
$ws.Range('E16').Value = '  +2.91%  '
$ws.Range('D17').Value = '60.25'
$ws.Range('E17').Value = '  +2.41%  '
$ws.Range('E18').Value = '  +2.03%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').Value = '201.27'
$ws.Range('E20').Value = '  +8.68%  '
$ws.Range('E21').Value = '  +2.60%  '
$ws.Range('D22').Value = '9.25'
$ws.Range('D23').Value = '5.99'
$ws.Range('E23').Value = '  +2.69%  '
$ws.Range('D24').Value = '1.83'
$ws.Range('E24').Value = '  +11.24%  '
$ws.Range('D25').Value = '141.37'
$ws.Range('E25').Value = '  +0.25%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('E27').Value = '  -5.24%  '
$ws.Range('E28').Value = '  +2.00%  '
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('E31').Value = '  +1.40%  '
$ws.Range('E32').Value = '  +2.50%  '
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('D34').Value = '1.47'
$ws.Range('E34').Value = '  +1.93%  '
$ws.Range('E35').Value = '  +1.14%  '
$ws.Range('E36').Value = '  +11.42%  '
$ws.Range('D37').Value = '1.126.13'
$ws.Range('E37').Value = '  +3.84%  '
$ws.Range('D39').Value = '0.792'
$ws.Range('E39').Value = '  +3.64%  '
$ws.Range('E40').Value = '  +2.41%  '
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('D42').Value = '0.781'
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('D43').Value = '1.735.61'
$ws.Range('D44').Value = '5.13'
$ws.Range('E44').Value = '  +1.57%  '
$ws.Range('D45').Value = '92.99'
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('E46').Value = '  +4.19%  '
$ws.Range('D47').Value = '53.29'
$ws.Range('E47').Value = '  +2.10%  '
$ws.Range('E48').Value = '  +0.19%  '
$ws.Range('D49').Value = '0.409'
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').Value = '0.0₇0925'
$ws.Range('E51').Value = '  -16.84%  '
